$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a literal string into a cell without letting Excel's COM
# automation auto-coerce date-shaped text (e.g. "2011-12-18") into a date
# serial number. We stage the text in a scratch cell far outside any used
# range, force it to Text format there, then Copy/PasteSpecial *values only*
# into the destination -- that leaves the destination's own style untouched
# (so it keeps matching whatever header/data style it already had).
# ---------------------------------------------------------------------------
function Set-TextValue($sheet, $row, $col, $text) {
    $helper = $sheet.Cells.Item(500, 500)
    $helper.NumberFormat = "@"
    $helper.Value = $text
    $helper.Copy() | Out-Null
    $sheet.Cells.Item($row, $col).PasteSpecial(-4163) | Out-Null
    $helper.Clear() | Out-Null
}

# ---------------------------------------------------------------------------
# Sheet "現金" (cash) -- 4th worksheet
# Add metadata columns (property_category, category, date, legislator_name,
# legislator_id, source_file, index) to the header row and to every data
# row, matching the layout already used on the other property sheets.
# ---------------------------------------------------------------------------
$wsCash = $wb.Worksheets.Item(4)

# Header row (row 1): copy the bold/border style from the existing D1
# header cell so the newly-created header cells pick up the same style
# (the engine dedups styles, so this lands on the same style index as the
# other header cells instead of creating new ones).
$wsCash.Cells.Item(1, 4).Copy($wsCash.Cells.Item(1, 5))
$wsCash.Cells.Item(1, 4).Copy($wsCash.Cells.Item(1, 6))
$wsCash.Cells.Item(1, 4).Copy($wsCash.Cells.Item(1, 7))
$wsCash.Cells.Item(1, 4).Copy($wsCash.Cells.Item(1, 8))
$wsCash.Cells.Item(1, 4).Copy($wsCash.Cells.Item(1, 9))
$wsCash.Cells.Item(1, 4).Copy($wsCash.Cells.Item(1, 10))
$wsCash.Cells.Item(1, 4).Copy($wsCash.Cells.Item(1, 11))

$wsCash.Cells.Item(1, 5).Value = "property_category"
$wsCash.Cells.Item(1, 6).Value = "category"
$wsCash.Cells.Item(1, 7).Value = "date"
$wsCash.Cells.Item(1, 8).Value = "legislator_name"
$wsCash.Cells.Item(1, 9).Value = "legislator_id"
$wsCash.Cells.Item(1, 10).Value = "source_file"
$wsCash.Cells.Item(1, 11).Value = "index"

# Row 2 (index 40)
$wsCash.Cells.Item(2, 5).Value = "cash"
$wsCash.Cells.Item(2, 6).Value = "normal"
Set-TextValue $wsCash 2 7 "2011-12-18"
$wsCash.Cells.Item(2, 8).Value = "林淑芬"
$wsCash.Cells.Item(2, 9).Value = 1337
$wsCash.Cells.Item(2, 10).Value = "tmp7b501"
$wsCash.Cells.Item(2, 11).Value = 40

# Row 3 (index 41) -- also fix D3, which used to store "1527000" as text;
# it becomes a real number. Owner stays 邱〇由.
$wsCash.Cells.Item(3, 4).Value = 1527000
$wsCash.Cells.Item(3, 5).Value = "cash"
$wsCash.Cells.Item(3, 6).Value = "normal"
Set-TextValue $wsCash 3 7 "2011-12-18"
$wsCash.Cells.Item(3, 8).Value = "林淑芬"
$wsCash.Cells.Item(3, 9).Value = 1337
$wsCash.Cells.Item(3, 10).Value = "tmp7b501"
$wsCash.Cells.Item(3, 11).Value = 41

# Row 4 (index 42)
$wsCash.Cells.Item(4, 5).Value = "cash"
$wsCash.Cells.Item(4, 6).Value = "normal"
Set-TextValue $wsCash 4 7 "2011-12-18"
$wsCash.Cells.Item(4, 8).Value = "林淑芬"
$wsCash.Cells.Item(4, 9).Value = 1337
$wsCash.Cells.Item(4, 10).Value = "tmp7b501"
$wsCash.Cells.Item(4, 11).Value = 42

# ---------------------------------------------------------------------------
# Sheet "存款" (bank deposit) -- 5th worksheet
# Same kind of metadata-column extension, plus real "bank" / "deposit_type"
# / "currency" header labels replacing the old placeholder header values.
# ---------------------------------------------------------------------------
$wsDeposit = $wb.Worksheets.Item(5)

$wsDeposit.Cells.Item(1, 4).Copy($wsDeposit.Cells.Item(1, 7))
$wsDeposit.Cells.Item(1, 4).Copy($wsDeposit.Cells.Item(1, 8))
$wsDeposit.Cells.Item(1, 4).Copy($wsDeposit.Cells.Item(1, 9))
$wsDeposit.Cells.Item(1, 4).Copy($wsDeposit.Cells.Item(1, 10))
$wsDeposit.Cells.Item(1, 4).Copy($wsDeposit.Cells.Item(1, 11))
$wsDeposit.Cells.Item(1, 4).Copy($wsDeposit.Cells.Item(1, 12))
$wsDeposit.Cells.Item(1, 4).Copy($wsDeposit.Cells.Item(1, 13))

$wsDeposit.Cells.Item(1, 2).Value = "bank"
$wsDeposit.Cells.Item(1, 3).Value = "deposit_type"
$wsDeposit.Cells.Item(1, 7).Value = "property_category"
$wsDeposit.Cells.Item(1, 8).Value = "category"
$wsDeposit.Cells.Item(1, 9).Value = "date"
$wsDeposit.Cells.Item(1, 10).Value = "legislator_name"
$wsDeposit.Cells.Item(1, 11).Value = "legislator_id"
$wsDeposit.Cells.Item(1, 12).Value = "source_file"
$wsDeposit.Cells.Item(1, 13).Value = "index"

# Row 2 (index 46)
$wsDeposit.Cells.Item(2, 7).Value = "deposit"
$wsDeposit.Cells.Item(2, 8).Value = "normal"
Set-TextValue $wsDeposit 2 9 "2011-12-18"
$wsDeposit.Cells.Item(2, 10).Value = "林淑芬"
$wsDeposit.Cells.Item(2, 11).Value = 1337
$wsDeposit.Cells.Item(2, 12).Value = "tmp7b501"
$wsDeposit.Cells.Item(2, 13).Value = 46

# Row 3 (index 47)
$wsDeposit.Cells.Item(3, 7).Value = "deposit"
$wsDeposit.Cells.Item(3, 8).Value = "normal"
Set-TextValue $wsDeposit 3 9 "2011-12-18"
$wsDeposit.Cells.Item(3, 10).Value = "林淑芬"
$wsDeposit.Cells.Item(3, 11).Value = 1337
$wsDeposit.Cells.Item(3, 12).Value = "tmp7b501"
$wsDeposit.Cells.Item(3, 13).Value = 47

# Row 4 (index 48)
$wsDeposit.Cells.Item(4, 7).Value = "deposit"
$wsDeposit.Cells.Item(4, 8).Value = "normal"
Set-TextValue $wsDeposit 4 9 "2011-12-18"
$wsDeposit.Cells.Item(4, 10).Value = "林淑芬"
$wsDeposit.Cells.Item(4, 11).Value = 1337
$wsDeposit.Cells.Item(4, 12).Value = "tmp7b501"
$wsDeposit.Cells.Item(4, 13).Value = 48
